# Append a new data row (row 6) to the tanker arrival log with the
# Potassium Chloride tanker entry. All cells in this sheet are stored as
# text (inline strings) in the original workbook, so values that look
# numeric/date-like ("0.12", "2026-02-17", "17:30", ...) must be forced
# to Text before assignment, otherwise Excel would auto-convert them to
# numbers / dates / times. Resetting the style back to "Normal" afterwards
# keeps the cell formatted as plain text while dropping the now-unneeded
# explicit number-format style, matching the rest of the sheet (which has
# no per-cell styling).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Text
    )
    $rng = $ws.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = "Normal"
}

Set-TextValue "A6" "GJ2010"
Set-TextValue "B6" "Potassium Chloride"
Set-TextValue "C6" "0.12"
Set-TextValue "D6" "abc"
Set-TextValue "E6" "djcsdj"
Set-TextValue "F6" "2026-02-17"
Set-TextValue "G6" "17:30"
Set-TextValue "H6" "2026-02-18"
Set-TextValue "I6" "17:30"
Set-TextValue "J6" "B101003"
Set-TextValue "K6" "O101"
Set-TextValue "L6" "Central Admin"
Set-TextValue "M6" "16-02-2026 17:27"
